# "add stunting OR given diarrhoea"
#
# 1. Rename "Incidence Diarrhoea" -> "Incidence diarrhoea" (lowercase "d").
# 2. Insert a new worksheet "OR stunting diarrhoea" right after
#    "RR diarrhoea" (and before "birth distribution"), holding the
#    age-band headers and a row of odds ratios (1.04 across the board).

$wb = $excel.ActiveWorkbook

# --- 1. Rename sheet ---------------------------------------------------
$incidence = $wb.Worksheets.Item("Incidence Diarrhoea")
$incidence.Name = "Incidence diarrhoea"

# --- 2. Insert the new sheet in the right slot -------------------------
$rrDiarrhoea = $wb.Worksheets.Item("RR diarrhoea")
$newSheet = $wb.Worksheets.Add($null, $rrDiarrhoea)
$newSheet.Name = "OR stunting diarrhoea"

$newSheet.Range("A1").Value = "<1 month"
$newSheet.Range("B1").Value = "1-5 months"
$newSheet.Range("C1").Value = "6-11 months"
$newSheet.Range("D1").Value = "12-23 months"
$newSheet.Range("E1").Value = "24-59 months"

$newSheet.Range("A2:E2").Value = 1.04
